{"js": "// Cover letter update:\n//  1. Date: \"September 09, 2024\" -> \"September 10, 2024\"\n//  2. First body paragraph: rewrite text\n//  3. Second body paragraph: rewrite text\n//  4. Third body paragraph: rewrite text\n//  5. Section margins: left/right 1800 -> 720 twips (1.25in -> 0.5in)\n\nconst replacements = [\n  {\n    oldText:\n      \"September 09, 2024\",\n    newText:\n      \"September 10, 2024\",\n  },\n  {\n    oldText:\n      \"I am excited to apply for the React Frontend Developer position at Rocketlane, where I can leverage my skills in building real-time collaboration experiences and crafting modern UI components. With a strong foundation in React, SPA Framework patterns, HTML, CSS, and SASS, I am confident in my ability to contribute to the development of innovative products. My experience in optimizing queries and using Redis to decrease loading time has allowed me to enhance user experience, and I am eager to bring this expertise to your team.\",\n    newText:\n      \"I am excited to apply for the React Frontend Developer position at Rocketlane, where I can leverage my skills in building real-time collaboration experiences and crafting modern UI components. With a strong foundation in React, SPA Framework patterns, HTML, CSS, and SASS, I am confident in my ability to contribute to the development of innovative products. My experience in building client-side applications with React.js and optimizing queries using Redis has allowed me to enhance user experience and improve application performance.\",\n  },\n  {\n    oldText:\n      \"As a detail-oriented and ambitious front-end engineer, I have developed a passion for creating flexible and interactive user interfaces. My experience as a Web Developer Intern at Digiidunia has given me hands-on experience in building client-side applications with React.js and CSS, and I am excited to apply this knowledge in a collaborative environment. Additionally, my experience in working with various technologies, including Express.js, MongoDB, and GitHub, has allowed me to develop a comprehensive understanding of the development process.\",\n    newText:\n      \"As a detail-oriented and ambitious frontend engineer, I am drawn to Rocketlane's mission to build real-time collaboration experiences. My experience as a Web Developer Intern at Digiidunia has provided me with hands-on experience in developing backend applications with Express.js and integrating MongoDB. Additionally, my personal projects, such as EcoSavvy, SmartLegalX, NutriSure, and PhishNet, demonstrate my ability to work with various technologies, including React.js, CSS, and MongoDB. I am excited about the opportunity to bring my skills and experience to Rocketlane and contribute to the development of cutting-edge products.\",\n  },\n  {\n    oldText:\n      \"I am particularly drawn to Rocketlane's commitment to building real-time collaboration experiences and its focus on innovation. As someone who is passionate about creating seamless and efficient user experiences, I believe that I would be a valuable addition to your team. I am impressed by Rocketlane's dedication to pushing the boundaries of what is possible and its commitment to excellence. Thank you for considering my application, and I look forward to the opportunity to contribute to the success of Rocketlane.\",\n    newText:\n      \"I am particularly impressed by Rocketlane's commitment to innovation and customer satisfaction. As someone who is passionate about building modern and efficient applications, I believe that I would be a great fit for the team. I am excited about the opportunity to work with a talented team of engineers and contribute to the development of products that make a real impact. Thank you for considering my application. I look forward to the opportunity to discuss my qualifications further.\",\n  },\n];\n\nconst body = context.document.body;\n\nfor (const { oldText, newText } of replacements) {\n  const results = body.search(oldText, { matchCase: true });\n  results.load(\"items\");\n  await context.sync();\n\n  if (results.items.length === 0) {\n    throw new Error(\"Could not find text to replace: \" + oldText.slice(0, 40));\n  }\n\n  results.items[0].insertText(newText, Word.InsertLocation.replace);\n}\n\nawait context.sync();\n\n// Section page margins: right/left go from 1800 -> 720 twips (0.5in = 720 twips @ 20 twips/pt, 36pt).\nconst sections = context.document.sections;\nsections.load(\"items\");\nawait context.sync();\n\nfor (const section of sections.items) {\n  section.pageSetup.leftMargin = 36; // points (720 twips)\n  section.pageSetup.rightMargin = 36; // points (720 twips)\n}\n\nawait context.sync();\n", "ps1": "# Cover letter update:\n#  1. Date: \"September 09, 2024\" -> \"September 10, 2024\"\n#  2. First body paragraph: rewrite text\n#  3. Second body paragraph: rewrite text\n#  4. Third body paragraph: rewrite text\n#  5. Section margins: left/right 1800 -> 720 twips (1.25in -> 0.5in, i.e. 90pt -> 36pt)\n\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @{\n        Old = \"September 09, 2024\"\n        New = \"September 10, 2024\"\n    },\n    @{\n        Old = \"I am excited to apply for the React Frontend Developer position at Rocketlane, where I can leverage my skills in building real-time collaboration experiences and crafting modern UI components. With a strong foundation in React, SPA Framework patterns, HTML, CSS, and SASS, I am confident in my ability to contribute to the development of innovative products. My experience in optimizing queries and using Redis to decrease loading time has allowed me to enhance user experience, and I am eager to bring this expertise to your team.\"\n        New = \"I am excited to apply for the React Frontend Developer position at Rocketlane, where I can leverage my skills in building real-time collaboration experiences and crafting modern UI components. With a strong foundation in React, SPA Framework patterns, HTML, CSS, and SASS, I am confident in my ability to contribute to the development of innovative products. My experience in building client-side applications with React.js and optimizing queries using Redis has allowed me to enhance user experience and improve application performance.\"\n    },\n    @{\n        Old = \"As a detail-oriented and ambitious front-end engineer, I have developed a passion for creating flexible and interactive user interfaces. My experience as a Web Developer Intern at Digiidunia has given me hands-on experience in building client-side applications with React.js and CSS, and I am excited to apply this knowledge in a collaborative environment. Additionally, my experience in working with various technologies, including Express.js, MongoDB, and GitHub, has allowed me to develop a comprehensive understanding of the development process.\"\n        New = \"As a detail-oriented and ambitious frontend engineer, I am drawn to Rocketlane's mission to build real-time collaboration experiences. My experience as a Web Developer Intern at Digiidunia has provided me with hands-on experience in developing backend applications with Express.js and integrating MongoDB. Additionally, my personal projects, such as EcoSavvy, SmartLegalX, NutriSure, and PhishNet, demonstrate my ability to work with various technologies, including React.js, CSS, and MongoDB. I am excited about the opportunity to bring my skills and experience to Rocketlane and contribute to the development of cutting-edge products.\"\n    },\n    @{\n        Old = \"I am particularly drawn to Rocketlane's commitment to building real-time collaboration experiences and its focus on innovation. As someone who is passionate about creating seamless and efficient user experiences, I believe that I would be a valuable addition to your team. I am impressed by Rocketlane's dedication to pushing the boundaries of what is possible and its commitment to excellence. Thank you for considering my application, and I look forward to the opportunity to contribute to the success of Rocketlane.\"\n        New = \"I am particularly impressed by Rocketlane's commitment to innovation and customer satisfaction. As someone who is passionate about building modern and efficient applications, I believe that I would be a great fit for the team. I am excited about the opportunity to work with a talented team of engineers and contribute to the development of products that make a real impact. Thank you for considering my application. I look forward to the opportunity to discuss my qualifications further.\"\n    }\n)\n\nforeach ($r in $replacements) {\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $result = $find.Execute($r.Old, $false, $false, $false, $false, $false, $true, 1, $false, $r.New, 2)\n    if (-not $result) {\n        throw \"Could not find text to replace: $($r.Old.Substring(0, 40))\"\n    }\n}\n\nforeach ($sec in $d.Sections) {\n    $sec.PageSetup.LeftMargin = 36\n    $sec.PageSetup.RightMargin = 36\n}\n"}
